$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Unit Processes")

# New row for the simplified EAF (Electric Arc Furnace) unit process,
# following the same pattern as the other "simple_*" steel units.
$ws.Range("A43").Value = "simple_EAF"
$ws.Range("B43").Value = "steel"
$ws.Range("C43").Value = "Electric Arc Furnace"
$ws.Range("D43").Value = "steel"
$ws.Range("E43").Value = "outflow"
$ws.Range("F43").Value = "data/steel/steel_simplified_var.xlsx"
$ws.Range("G43").Value = "EAF"
$ws.Range("H43").Value = "data/steel/steel_simplified_calcs.xlsx"
$ws.Range("I43").Value = "EAF"

$ws.Range("A43").NumberFormat = "@"
$ws.Range("C43:I43").NumberFormat = "@"

$ws.Range("I43").Select()
